# Auto-generated Excel COM-interop edit script
# Applies the "plotted toy-spam with min 5" update to avg_0.004_scores.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update every cell value in the two result tables (A:H and J:Q) ---
$ws.Cells.Item(1,1).Value = 'negative'
$ws.Cells.Item(1,10).Value = 'positive'
$ws.Cells.Item(2,1).Value = 'name'
$ws.Cells.Item(2,2).Value = 'anchor score'
$ws.Cells.Item(2,3).Value = 'type occurences'
$ws.Cells.Item(2,4).Value = 'total occurences'
$ws.Cells.Item(2,5).Value = '+%'
$ws.Cells.Item(2,6).Value = '-%'
$ws.Cells.Item(2,7).Value = 'both'
$ws.Cells.Item(2,8).Value = 'normal'
$ws.Cells.Item(2,10).Value = 'name'
$ws.Cells.Item(2,11).Value = 'anchor score'
$ws.Cells.Item(2,12).Value = 'type occurences'
$ws.Cells.Item(2,13).Value = 'total occurences'
$ws.Cells.Item(2,14).Value = '+%'
$ws.Cells.Item(2,15).Value = '-%'
$ws.Cells.Item(2,16).Value = 'both'
$ws.Cells.Item(2,17).Value = 'normal'
$ws.Cells.Item(3,1).Value = 'poorly'
$ws.Cells.Item(3,2).Value = 0.9565217391304348
$ws.Cells.Item(3,3).Value = 44.0
$ws.Cells.Item(3,4).Value = 44.0
$ws.Cells.Item(3,5).Value = 0.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = $false
$ws.Cells.Item(3,8).Value = 2.0
$ws.Cells.Item(3,10).Value = 'wonderful'
$ws.Cells.Item(3,11).Value = 0.875
$ws.Cells.Item(3,12).Value = 49.0
$ws.Cells.Item(3,13).Value = 49.0
$ws.Cells.Item(3,14).Value = 1.0
$ws.Cells.Item(3,15).Value = 0.0
$ws.Cells.Item(3,16).Value = $false
$ws.Cells.Item(3,17).Value = 7.0
$ws.Cells.Item(4,1).Value = 'disappointing'
$ws.Cells.Item(4,2).Value = 0.7954545454545454
$ws.Cells.Item(4,3).Value = 35.0
$ws.Cells.Item(4,4).Value = 35.0
$ws.Cells.Item(4,5).Value = 0.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = $false
$ws.Cells.Item(4,8).Value = 9.0
$ws.Cells.Item(4,10).Value = 'awesome'
$ws.Cells.Item(4,11).Value = 0.8153846153846154
$ws.Cells.Item(4,12).Value = 53.0
$ws.Cells.Item(4,13).Value = 53.0
$ws.Cells.Item(4,14).Value = 1.0
$ws.Cells.Item(4,15).Value = 0.0
$ws.Cells.Item(4,16).Value = $false
$ws.Cells.Item(4,17).Value = 12.0
$ws.Cells.Item(5,1).Value = 'returned'
$ws.Cells.Item(5,2).Value = 0.7894736842105263
$ws.Cells.Item(5,3).Value = 30.0
$ws.Cells.Item(5,4).Value = 30.0
$ws.Cells.Item(5,5).Value = 0.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = $false
$ws.Cells.Item(5,8).Value = 8.0
$ws.Cells.Item(5,10).Value = 'favorite'
$ws.Cells.Item(5,11).Value = 0.7204301075268817
$ws.Cells.Item(5,12).Value = 67.0
$ws.Cells.Item(5,13).Value = 67.0
$ws.Cells.Item(5,14).Value = 1.0
$ws.Cells.Item(5,15).Value = 0.0
$ws.Cells.Item(5,16).Value = $false
$ws.Cells.Item(5,17).Value = 26.0
$ws.Cells.Item(6,1).Value = 'however'
$ws.Cells.Item(6,2).Value = 0.75
$ws.Cells.Item(6,3).Value = 48.0
$ws.Cells.Item(6,4).Value = 48.0
$ws.Cells.Item(6,5).Value = 0.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = $false
$ws.Cells.Item(6,8).Value = 16.0
$ws.Cells.Item(6,10).Value = 'classic'
$ws.Cells.Item(6,11).Value = 0.6415094339622641
$ws.Cells.Item(6,12).Value = 34.0
$ws.Cells.Item(6,13).Value = 34.0
$ws.Cells.Item(6,14).Value = 1.0
$ws.Cells.Item(6,15).Value = 0.0
$ws.Cells.Item(6,16).Value = $false
$ws.Cells.Item(6,17).Value = 19.0
$ws.Cells.Item(7,1).Value = 'poor'
$ws.Cells.Item(7,2).Value = 0.7464788732394366
$ws.Cells.Item(7,3).Value = 53.0
$ws.Cells.Item(7,4).Value = 53.0
$ws.Cells.Item(7,5).Value = 0.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = $false
$ws.Cells.Item(7,8).Value = 18.0
$ws.Cells.Item(7,10).Value = 'excellent'
$ws.Cells.Item(7,11).Value = 0.484375
$ws.Cells.Item(7,12).Value = 31.0
$ws.Cells.Item(7,13).Value = 31.0
$ws.Cells.Item(7,14).Value = 1.0
$ws.Cells.Item(7,15).Value = 0.0
$ws.Cells.Item(7,16).Value = $false
$ws.Cells.Item(7,17).Value = 33.0
$ws.Cells.Item(8,1).Value = 'disappointed'
$ws.Cells.Item(8,2).Value = 0.7419354838709677
$ws.Cells.Item(8,3).Value = 138.0
$ws.Cells.Item(8,4).Value = 138.0
$ws.Cells.Item(8,5).Value = 0.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = $false
$ws.Cells.Item(8,8).Value = 48.0
$ws.Cells.Item(8,10).Value = 'thank'
$ws.Cells.Item(8,11).Value = 0.4782608695652174
$ws.Cells.Item(8,12).Value = 33.0
$ws.Cells.Item(8,13).Value = 33.0
$ws.Cells.Item(8,14).Value = 1.0
$ws.Cells.Item(8,15).Value = 0.0
$ws.Cells.Item(8,16).Value = $false
$ws.Cells.Item(8,17).Value = 36.0
$ws.Cells.Item(9,1).Value = 'broke'
$ws.Cells.Item(9,2).Value = 0.7038834951456311
$ws.Cells.Item(9,3).Value = 145.0
$ws.Cells.Item(9,4).Value = 145.0
$ws.Cells.Item(9,5).Value = 0.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = $false
$ws.Cells.Item(9,8).Value = 61.0
$ws.Cells.Item(9,10).Value = 'great'
$ws.Cells.Item(9,11).Value = 0.3795081967213115
$ws.Cells.Item(9,12).Value = 463.0
$ws.Cells.Item(9,13).Value = 463.0
$ws.Cells.Item(9,14).Value = 1.0
$ws.Cells.Item(9,15).Value = 0.0
$ws.Cells.Item(9,16).Value = $false
$ws.Cells.Item(9,17).Value = 757.0
$ws.Cells.Item(10,1).Value = 'instead'
$ws.Cells.Item(10,2).Value = 0.6458333333333334
$ws.Cells.Item(10,3).Value = 31.0
$ws.Cells.Item(10,4).Value = 31.0
$ws.Cells.Item(10,5).Value = 0.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = $false
$ws.Cells.Item(10,8).Value = 17.0
$ws.Cells.Item(10,10).Value = 'love'
$ws.Cells.Item(10,11).Value = 0.3242467718794835
$ws.Cells.Item(10,12).Value = 226.0
$ws.Cells.Item(10,13).Value = 226.0
$ws.Cells.Item(10,14).Value = 1.0
$ws.Cells.Item(10,15).Value = 0.0
$ws.Cells.Item(10,16).Value = $false
$ws.Cells.Item(10,17).Value = 471.0
$ws.Cells.Item(11,1).Value = 'waste'
$ws.Cells.Item(11,2).Value = 0.6418918918918919
$ws.Cells.Item(11,3).Value = 95.0
$ws.Cells.Item(11,4).Value = 95.0
$ws.Cells.Item(11,5).Value = 0.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = $false
$ws.Cells.Item(11,8).Value = 53.0
$ws.Cells.Item(11,10).Value = 'loves'
$ws.Cells.Item(11,11).Value = 0.2904564315352697
$ws.Cells.Item(11,12).Value = 140.0
$ws.Cells.Item(11,13).Value = 140.0
$ws.Cells.Item(11,14).Value = 1.0
$ws.Cells.Item(11,15).Value = 0.0
$ws.Cells.Item(11,16).Value = $false
$ws.Cells.Item(11,17).Value = 342.0
$ws.Cells.Item(12,1).Value = 'smaller'
$ws.Cells.Item(12,2).Value = 0.6134453781512605
$ws.Cells.Item(12,3).Value = 73.0
$ws.Cells.Item(12,4).Value = 73.0
$ws.Cells.Item(12,5).Value = 0.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = $false
$ws.Cells.Item(12,8).Value = 46.0
$ws.Cells.Item(12,10).Value = 'loved'
$ws.Cells.Item(12,11).Value = 0.2293577981651376
$ws.Cells.Item(12,12).Value = 75.0
$ws.Cells.Item(12,13).Value = 75.0
$ws.Cells.Item(12,14).Value = 1.0
$ws.Cells.Item(12,15).Value = 0.0
$ws.Cells.Item(12,16).Value = $false
$ws.Cells.Item(12,17).Value = 252.0
$ws.Cells.Item(13,1).Value = 'junk'
$ws.Cells.Item(13,2).Value = 0.5818181818181818
$ws.Cells.Item(13,3).Value = 32.0
$ws.Cells.Item(13,4).Value = 32.0
$ws.Cells.Item(13,5).Value = 0.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = $false
$ws.Cells.Item(13,8).Value = 23.0
$ws.Cells.Item(13,10).Value = 'perfect'
$ws.Cells.Item(13,11).Value = 0.2168674698795181
$ws.Cells.Item(13,12).Value = 36.0
$ws.Cells.Item(13,13).Value = 36.0
$ws.Cells.Item(13,14).Value = 1.0
$ws.Cells.Item(13,15).Value = 0.0
$ws.Cells.Item(13,16).Value = $false
$ws.Cells.Item(13,17).Value = 130.0
$ws.Cells.Item(14,1).Value = 'small'
$ws.Cells.Item(14,2).Value = 0.5101449275362319
$ws.Cells.Item(14,3).Value = 176.0
$ws.Cells.Item(14,4).Value = 176.0
$ws.Cells.Item(14,5).Value = 0.0
$ws.Cells.Item(14,6).Value = 1.0
$ws.Cells.Item(14,7).Value = $false
$ws.Cells.Item(14,8).Value = 169.0
$ws.Cells.Item(14,10).Value = 'christmas'
$ws.Cells.Item(14,11).Value = 0.1204819277108434
$ws.Cells.Item(14,12).Value = 30.0
$ws.Cells.Item(14,13).Value = 30.0
$ws.Cells.Item(14,14).Value = 1.0
$ws.Cells.Item(14,15).Value = 0.0
$ws.Cells.Item(14,16).Value = $false
$ws.Cells.Item(14,17).Value = 219.0
$ws.Cells.Item(15,1).Value = 'broken'
$ws.Cells.Item(15,2).Value = 0.4819277108433735
$ws.Cells.Item(15,3).Value = 40.0
$ws.Cells.Item(15,4).Value = 40.0
$ws.Cells.Item(15,5).Value = 0.0
$ws.Cells.Item(15,6).Value = 1.0
$ws.Cells.Item(15,7).Value = $false
$ws.Cells.Item(15,8).Value = 43.0
$ws.Cells.Item(15,10).Value = 'fun'
$ws.Cells.Item(15,11).Value = 0.1130587204206836
$ws.Cells.Item(15,12).Value = 129.0
$ws.Cells.Item(15,13).Value = 129.0
$ws.Cells.Item(15,14).Value = 1.0
$ws.Cells.Item(15,15).Value = 0.0
$ws.Cells.Item(15,16).Value = $false
$ws.Cells.Item(15,17).Value = 1012.0
$ws.Cells.Item(16,1).Value = 'paint'
$ws.Cells.Item(16,2).Value = 0.4603174603174603
$ws.Cells.Item(16,3).Value = 29.0
$ws.Cells.Item(16,4).Value = 29.0
$ws.Cells.Item(16,5).Value = 0.0
$ws.Cells.Item(16,6).Value = 1.0
$ws.Cells.Item(16,7).Value = $false
$ws.Cells.Item(16,8).Value = 34.0
$ws.Cells.Item(16,10).Value = 'game'
$ws.Cells.Item(16,11).Value = 0.05191434133679429
$ws.Cells.Item(16,12).Value = 80.0
$ws.Cells.Item(16,13).Value = 80.0
$ws.Cells.Item(16,14).Value = 1.0
$ws.Cells.Item(16,15).Value = 0.0
$ws.Cells.Item(16,16).Value = $false
$ws.Cells.Item(16,17).Value = 1461.0
$ws.Cells.Item(17,1).Value = 'di'
$ws.Cells.Item(17,2).Value = 0.453125
$ws.Cells.Item(17,3).Value = 29.0
$ws.Cells.Item(17,4).Value = 29.0
$ws.Cells.Item(17,5).Value = 0.0
$ws.Cells.Item(17,6).Value = 1.0
$ws.Cells.Item(17,7).Value = $false
$ws.Cells.Item(17,8).Value = 35.0
$ws.Cells.Item(18,1).Value = 'apart'
$ws.Cells.Item(18,2).Value = 0.4421052631578947
$ws.Cells.Item(18,3).Value = 42.0
$ws.Cells.Item(18,4).Value = 42.0
$ws.Cells.Item(18,5).Value = 0.0
$ws.Cells.Item(18,6).Value = 1.0
$ws.Cells.Item(18,7).Value = $false
$ws.Cells.Item(18,8).Value = 53.0
$ws.Cells.Item(19,1).Value = 'plastic'
$ws.Cells.Item(19,2).Value = 0.4330708661417323
$ws.Cells.Item(19,3).Value = 55.0
$ws.Cells.Item(19,4).Value = 55.0
$ws.Cells.Item(19,5).Value = 0.0
$ws.Cells.Item(19,6).Value = 1.0
$ws.Cells.Item(19,7).Value = $false
$ws.Cells.Item(19,8).Value = 72.0
$ws.Cells.Item(20,1).Value = 'difficult'
$ws.Cells.Item(20,2).Value = 0.4044943820224719
$ws.Cells.Item(20,3).Value = 36.0
$ws.Cells.Item(20,4).Value = 36.0
$ws.Cells.Item(20,5).Value = 0.0
$ws.Cells.Item(20,6).Value = 1.0
$ws.Cells.Item(20,7).Value = $false
$ws.Cells.Item(20,8).Value = 53.0
$ws.Cells.Item(21,1).Value = 'ok'
$ws.Cells.Item(21,2).Value = 0.3984375
$ws.Cells.Item(21,3).Value = 51.0
$ws.Cells.Item(21,4).Value = 51.0
$ws.Cells.Item(21,5).Value = 0.0
$ws.Cells.Item(21,6).Value = 1.0
$ws.Cells.Item(21,7).Value = $false
$ws.Cells.Item(21,8).Value = 77.0
$ws.Cells.Item(22,1).Value = 'thought'
$ws.Cells.Item(22,2).Value = 0.3316831683168317
$ws.Cells.Item(22,3).Value = 67.0
$ws.Cells.Item(22,4).Value = 67.0
$ws.Cells.Item(22,5).Value = 0.0
$ws.Cells.Item(22,6).Value = 1.0
$ws.Cells.Item(22,7).Value = $false
$ws.Cells.Item(22,8).Value = 135.0
$ws.Cells.Item(23,1).Value = 'cheap'
$ws.Cells.Item(23,2).Value = 0.3270142180094787
$ws.Cells.Item(23,3).Value = 69.0
$ws.Cells.Item(23,4).Value = 69.0
$ws.Cells.Item(23,5).Value = 0.0
$ws.Cells.Item(23,6).Value = 1.0
$ws.Cells.Item(23,7).Value = $false
$ws.Cells.Item(23,8).Value = 142.0
$ws.Cells.Item(24,1).Value = 'though'
$ws.Cells.Item(24,2).Value = 0.2991452991452991
$ws.Cells.Item(24,3).Value = 35.0
$ws.Cells.Item(24,4).Value = 35.0
$ws.Cells.Item(24,5).Value = 0.0
$ws.Cells.Item(24,6).Value = 1.0
$ws.Cells.Item(24,7).Value = $false
$ws.Cells.Item(24,8).Value = 82.0
$ws.Cells.Item(25,1).Value = 'size'
$ws.Cells.Item(25,2).Value = 0.2525773195876289
$ws.Cells.Item(25,3).Value = 49.0
$ws.Cells.Item(25,4).Value = 49.0
$ws.Cells.Item(25,5).Value = 0.0
$ws.Cells.Item(25,6).Value = 1.0
$ws.Cells.Item(25,7).Value = $false
$ws.Cells.Item(25,8).Value = 145.0
$ws.Cells.Item(26,1).Value = 'could'
$ws.Cells.Item(26,2).Value = 0.2229299363057325
$ws.Cells.Item(26,3).Value = 35.0
$ws.Cells.Item(26,4).Value = 35.0
$ws.Cells.Item(26,5).Value = 0.0
$ws.Cells.Item(26,6).Value = 1.0
$ws.Cells.Item(26,7).Value = $false
$ws.Cells.Item(26,8).Value = 122.0
$ws.Cells.Item(27,1).Value = 'money'
$ws.Cells.Item(27,2).Value = 0.2183544303797468
$ws.Cells.Item(27,3).Value = 69.0
$ws.Cells.Item(27,4).Value = 69.0
$ws.Cells.Item(27,5).Value = 0.0
$ws.Cells.Item(27,6).Value = 1.0
$ws.Cells.Item(27,7).Value = $false
$ws.Cells.Item(27,8).Value = 247.0
$ws.Cells.Item(28,1).Value = 'would'
$ws.Cells.Item(28,2).Value = 0.2106824925816024
$ws.Cells.Item(28,3).Value = 142.0
$ws.Cells.Item(28,4).Value = 142.0
$ws.Cells.Item(28,5).Value = 0.0
$ws.Cells.Item(28,6).Value = 1.0
$ws.Cells.Item(28,7).Value = $false
$ws.Cells.Item(28,8).Value = 532.0
$ws.Cells.Item(29,1).Value = 'used'
$ws.Cells.Item(29,2).Value = 0.1885714285714286
$ws.Cells.Item(29,3).Value = 33.0
$ws.Cells.Item(29,4).Value = 33.0
$ws.Cells.Item(29,5).Value = 0.0
$ws.Cells.Item(29,6).Value = 1.0
$ws.Cells.Item(29,7).Value = $false
$ws.Cells.Item(29,8).Value = 142.0
$ws.Cells.Item(30,1).Value = 'better'
$ws.Cells.Item(30,2).Value = 0.1869158878504673
$ws.Cells.Item(30,3).Value = 40.0
$ws.Cells.Item(30,4).Value = 40.0
$ws.Cells.Item(30,5).Value = 0.0
$ws.Cells.Item(30,6).Value = 1.0
$ws.Cells.Item(30,7).Value = $false
$ws.Cells.Item(30,8).Value = 174.0
$ws.Cells.Item(31,1).Value = 'hard'
$ws.Cells.Item(31,2).Value = 0.185
$ws.Cells.Item(31,3).Value = 37.0
$ws.Cells.Item(31,4).Value = 37.0
$ws.Cells.Item(31,5).Value = 0.0
$ws.Cells.Item(31,6).Value = 1.0
$ws.Cells.Item(31,7).Value = $false
$ws.Cells.Item(31,8).Value = 163.0
$ws.Cells.Item(32,1).Value = 'work'
$ws.Cells.Item(32,2).Value = 0.180379746835443
$ws.Cells.Item(32,3).Value = 57.0
$ws.Cells.Item(32,4).Value = 57.0
$ws.Cells.Item(32,5).Value = 0.0
$ws.Cells.Item(32,6).Value = 1.0
$ws.Cells.Item(32,7).Value = $false
$ws.Cells.Item(32,8).Value = 259.0
$ws.Cells.Item(33,1).Value = 'product'
$ws.Cells.Item(33,2).Value = 0.1784140969162996
$ws.Cells.Item(33,3).Value = 81.0
$ws.Cells.Item(33,4).Value = 81.0
$ws.Cells.Item(33,5).Value = 0.0
$ws.Cells.Item(33,6).Value = 1.0
$ws.Cells.Item(33,7).Value = $false
$ws.Cells.Item(33,8).Value = 373.0
$ws.Cells.Item(34,1).Value = 'item'
$ws.Cells.Item(34,2).Value = 0.1666666666666667
$ws.Cells.Item(34,3).Value = 46.0
$ws.Cells.Item(34,4).Value = 46.0
$ws.Cells.Item(34,5).Value = 0.0
$ws.Cells.Item(34,6).Value = 1.0
$ws.Cells.Item(34,7).Value = $false
$ws.Cells.Item(34,8).Value = 230.0
$ws.Cells.Item(35,1).Value = 'price'
$ws.Cells.Item(35,2).Value = 0.1354466858789625
$ws.Cells.Item(35,3).Value = 47.0
$ws.Cells.Item(35,4).Value = 48.0
$ws.Cells.Item(35,5).Value = 0.02
$ws.Cells.Item(35,6).Value = 0.98
$ws.Cells.Item(35,7).Value = $true
$ws.Cells.Item(35,8).Value = 300.0
$ws.Cells.Item(36,1).Value = '2'
$ws.Cells.Item(36,2).Value = 0.1278195488721804
$ws.Cells.Item(36,3).Value = 34.0
$ws.Cells.Item(36,4).Value = 35.0
$ws.Cells.Item(36,5).Value = 0.03
$ws.Cells.Item(36,6).Value = 0.97
$ws.Cells.Item(36,7).Value = $true
$ws.Cells.Item(36,8).Value = 232.0
$ws.Cells.Item(37,1).Value = '3'
$ws.Cells.Item(37,2).Value = 0.1255060728744939
$ws.Cells.Item(37,3).Value = 31.0
$ws.Cells.Item(37,4).Value = 32.0
$ws.Cells.Item(37,5).Value = 0.03
$ws.Cells.Item(37,6).Value = 0.97
$ws.Cells.Item(37,7).Value = $true
$ws.Cells.Item(37,8).Value = 216.0
$ws.Cells.Item(38,1).Value = 'use'
$ws.Cells.Item(38,2).Value = 0.1095890410958904
$ws.Cells.Item(38,3).Value = 40.0
$ws.Cells.Item(38,4).Value = 40.0
$ws.Cells.Item(38,5).Value = 0.0
$ws.Cells.Item(38,6).Value = 1.0
$ws.Cells.Item(38,7).Value = $false
$ws.Cells.Item(38,8).Value = 325.0
$ws.Cells.Item(39,1).Value = 'buy'
$ws.Cells.Item(39,2).Value = 0.08169014084507042
$ws.Cells.Item(39,3).Value = 29.0
$ws.Cells.Item(39,4).Value = 29.0
$ws.Cells.Item(39,5).Value = 0.0
$ws.Cells.Item(39,6).Value = 1.0
$ws.Cells.Item(39,7).Value = $false
$ws.Cells.Item(39,8).Value = 326.0
$ws.Cells.Item(40,1).Value = 'like'
$ws.Cells.Item(40,2).Value = 0.08072487644151564
$ws.Cells.Item(40,3).Value = 49.0
$ws.Cells.Item(40,4).Value = 50.0
$ws.Cells.Item(40,5).Value = 0.02
$ws.Cells.Item(40,6).Value = 0.98
$ws.Cells.Item(40,7).Value = $true
$ws.Cells.Item(40,8).Value = 558.0
$ws.Cells.Item(41,1).Value = 'little'
$ws.Cells.Item(41,2).Value = 0.06904231625835189
$ws.Cells.Item(41,3).Value = 31.0
$ws.Cells.Item(41,4).Value = 31.0
$ws.Cells.Item(41,5).Value = 0.0
$ws.Cells.Item(41,6).Value = 1.0
$ws.Cells.Item(41,7).Value = $false
$ws.Cells.Item(41,8).Value = 418.0
$ws.Cells.Item(42,1).Value = 'one'
$ws.Cells.Item(42,2).Value = 0.03929024081115336
$ws.Cells.Item(42,3).Value = 31.0
$ws.Cells.Item(42,4).Value = 36.0
$ws.Cells.Item(42,5).Value = 0.14
$ws.Cells.Item(42,6).Value = 0.86
$ws.Cells.Item(42,7).Value = $true
$ws.Cells.Item(42,8).Value = 758.0

# --- Apply header formatting (bold font + border + centered alignment) to the newly
# added header-style cells, by copying the format from an existing styled cell ---
$ws.Range("A3").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("J3").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

